$d = $word.ActiveDocument

# 1. Title / heading text (appears twice: H1 heading and bold teaser near the end)
$d.Content.Find.Execute(
    "Play Columbus for Free - Historical Slot Game", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Columbus Slot Free - Historical Theme with Wilds and Bonus Game", 2) | Out-Null

# 2. "What we like" bullet list
$d.Content.Find.Execute(
    "Historical theme with engaging symbols", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Historical theme", 2) | Out-Null

$d.Content.Find.Execute(
    "Wild symbol and Scatter symbol increase winning potential", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Wild symbol, Scatter symbol, and Bonus Game", 2) | Out-Null

$d.Content.Find.Execute(
    "Bonus Game offers up to ten free spins", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Autoplay function", 2) | Out-Null

$d.Content.Find.Execute(
    "Versatile betting options for different budgets", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Good winning potential", 2) | Out-Null

# 3. "What we don't like" bullet list
$d.Content.Find.Execute(
    "Only 9 paylines may limit gameplay", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Limited number of paylines", 2) | Out-Null

$d.Content.Find.Execute(
    "Graphics and sound effects could be improved", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Only one Wild symbol", 2) | Out-Null

# 4. Meta description (italic paragraph)
$d.Content.Find.Execute(
    "Read a review of the historical slot game Columbus. Play for free and enjoy a Bonus Game with up to ten free spins and a 20,000 coin jackpot.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Read our review of Columbus slot game with historical theme, Wilds, and Bonus Game. Play for free now!", 2) | Out-Null
